# Rebuild the LR-pairs data table (rows 2-17) with the corrected per-cluster-pair
# values supplied by Dr Hou (natmi re-run): every (sending cluster) x (target cluster)
# combination now appears exactly once, in ECs/FAPs/M2/sCs x ECs/FAPs/M2/sCs order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Jag1", "Notch3", "ECs", 3, 1, 16.12468433333333, 48.374053, 0.2955490655206278, 0.2955490655206279, 3, 1, 17.32492333333333, 51.97477, 0.1811369435820938, 0.1811369435820938, 279.3589198492011, 2514.23027864281, 0.05353485440695051, 0.05353485440695053),
    @("ECs", "Jag1", "Notch3", "FAPs", 3, 1, 16.12468433333333, 48.374053, 0.2955490655206278, 0.2955490655206279, 3, 1, 6.590479666666667, 19.771439, 0.06890531753540824, 0.06890531753540824, 106.2694042302519, 956.4246380722669, 0.02036490220699204, 0.02036490220699204),
    @("ECs", "Jag1", "Notch3", "M2", 3, 1, 16.12468433333333, 48.374053, 0.2955490655206278, 0.2955490655206279, 2, 0.6666666666666666, 0.111195, 0.333585, 0.001162574982531578, 0.001162574982531578, 1.792984274445, 16.136858470005, 0.000343597949684868, 0.0003435979496848681),
    @("ECs", "Jag1", "Notch3", "sCs", 3, 1, 16.12468433333333, 48.374053, 0.2955490655206278, 0.2955490655206279, 3, 1, 71.61884566666667, 214.856537, 0.7487951638999665, 0.7487951638999664, 1154.831278692718, 10393.48150823446, 0.2213057109570004, 0.2213057109570004),
    @("FAPs", "Jag1", "Notch3", "ECs", 3, 1, 14.68975, 44.06925, 0.2692481784748309, 0.2692481784748309, 3, 1, 17.32492333333333, 51.97477, 0.1811369435820938, 0.1811369435820938, 254.4987925358333, 2290.4891328225, 0.04877079211397698, 0.04877079211397698),
    @("FAPs", "Jag1", "Notch3", "FAPs", 3, 1, 14.68975, 44.06925, 0.2692481784748309, 0.2692481784748309, 3, 1, 6.590479666666667, 19.771439, 0.06890531753540824, 0.06890531753540824, 96.81249868341666, 871.31248815075, 0.01855263123363849, 0.01855263123363849),
    @("FAPs", "Jag1", "Notch3", "M2", 3, 1, 14.68975, 44.06925, 0.2692481784748309, 0.2692481784748309, 2, 0.6666666666666666, 0.111195, 0.333585, 0.001162574982531578, 0.001162574982531578, 1.63342675125, 14.70084076125, 0.0003130211963870357, 0.0003130211963870357),
    @("FAPs", "Jag1", "Notch3", "sCs", 3, 1, 14.68975, 44.06925, 0.2692481784748309, 0.2692481784748309, 3, 1, 71.61884566666667, 214.856537, 0.7487951638999665, 0.7487951638999664, 1052.062938131917, 9468.56644318725, 0.2016117339308284, 0.2016117339308284),
    @("M2", "Jag1", "Notch3", "ECs", 3, 1, 1.967860666666667, 5.903582, 0.03606888476606249, 0.03606888476606249, 3, 1, 17.32492333333333, 51.97477, 0.1811369435820938, 0.1811369435820938, 34.09303518068222, 306.83731662614, 0.006533407544939304, 0.006533407544939304),
    @("M2", "Jag1", "Notch3", "FAPs", 3, 1, 1.967860666666667, 5.903582, 0.03606888476606249, 0.03606888476606249, 3, 1, 6.590479666666667, 19.771439, 0.06890531753540824, 0.06890531753540824, 12.96914571049978, 116.722311394498, 0.002485337957953585, 0.002485337957953585),
    @("M2", "Jag1", "Notch3", "M2", 3, 1, 1.967860666666667, 5.903582, 0.03606888476606249, 0.03606888476606249, 2, 0.6666666666666666, 0.111195, 0.333585, 0.001162574982531578, 0.001162574982531578, 0.21881626683, 1.96934640147, 0.00004193278307683859, 0.00004193278307683859),
    @("M2", "Jag1", "Notch3", "sCs", 3, 1, 1.967860666666667, 5.903582, 0.03606888476606249, 0.03606888476606249, 3, 1, 71.61884566666667, 214.856537, 0.7487951638999665, 0.7487951638999664, 140.9359093795038, 1268.423184415534, 0.02700820648009276, 0.02700820648009276),
    @("sCs", "Jag1", "Notch3", "ECs", 3, 1, 21.77610566666667, 65.328317, 0.3991338712384788, 0.3991338712384788, 3, 1, 17.32492333333333, 51.97477, 0.1811369435820938, 0.1811369435820938, 377.2693611735656, 3395.42425056209, 0.07229788951622704, 0.07229788951622704),
    @("sCs", "Jag1", "Notch3", "FAPs", 3, 1, 21.77610566666667, 65.328317, 0.3991338712384788, 0.3991338712384788, 3, 1, 6.590479666666667, 19.771439, 0.06890531753540824, 0.06890531753540824, 143.5149816153514, 1291.634834538163, 0.02750244613682412, 0.02750244613682412),
    @("sCs", "Jag1", "Notch3", "M2", 3, 1, 21.77610566666667, 65.328317, 0.3991338712384788, 0.3991338712384788, 2, 0.6666666666666666, 0.111195, 0.333585, 0.001162574982531578, 0.001162574982531578, 2.421394069605, 21.792546626445, 0.0004640230533828355, 0.0004640230533828355),
    @("sCs", "Jag1", "Notch3", "sCs", 3, 1, 21.77610566666667, 65.328317, 0.3991338712384788, 0.3991338712384788, 3, 1, 71.61884566666667, 214.856537, 0.7487951638999665, 0.7487951638999664, 1559.579550962025, 14036.21595865823, 0.2988695125320448, 0.2988695125320448)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}
